$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: replace 机械臂1..5 labels (N1,P1,R1,T1,V1) with new arm-segment labels
$ws.Range("R1").Value = "中臂"
$ws.Range("T1").Value = "小臂"
$ws.Range("V1").Value = "夹手"
$ws.Range("P1").Value = "大臂"
$ws.Range("N1").Value = "水平"

# New data cell added at J12 (raw byte-frame sample string)
$ws.Range("J12").Value = "25 04 D0 05 DC 05 DC 05 DC 05 DC 05 DC 05 DC 05 DC 05 DC 05 DC 05 DC 05 DC 05 DC 08 00 21"

# Update view: scroll right so column I is the leftmost visible column, and move the
# active selection to O21
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("O21").Select()
